$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "41.670.32"
Set-TextValue "E2" "  -0.81%  "
Set-TextValue "D3" "2.222.62"
Set-TextValue "E3" "  -0.50%  "
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "250.94"
Set-TextValue "E5" "  +8.18%  "
Set-TextValue "D6" "0.628"
Set-TextValue "E6" "  -0.29%  "
Set-TextValue "D7" "70.94"
Set-TextValue "E7" "  +3.44%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "D9" "0.578"
Set-TextValue "E9" "  +5.19%  "
Set-TextValue "D10" "41.36"
Set-TextValue "E10" "  +18.30%  "
Set-TextValue "D11" "0.0965"
Set-TextValue "E11" "  -1.72%  "
Set-TextValue "D12" "58.19"
Set-TextValue "E12" "  +0.38%  "
Set-TextValue "E13" "  +0.40%  "
Set-TextValue "D14" "6.99"
Set-TextValue "E14" "  +4.17%  "
Set-TextValue "D15" "2.552.81"
Set-TextValue "E15" "  -0.78%  "
Set-TextValue "D16" "14.88"
Set-TextValue "E16" "  -0.04%  "
Set-TextValue "D17" "0.855"
Set-TextValue "E17" "  +0.37%  "
Set-TextValue "D18" "2.222.46"
Set-TextValue "E18" "  -0.54%  "
Set-TextValue "D19" "41.568.82"
Set-TextValue "E19" "  -0.57%  "
Set-TextValue "D20" "0.0₃0967"
Set-TextValue "E20" "  +0.09%  "
Set-TextValue "D21" "6.20"
Set-TextValue "E21" "  -0.10%  "
Set-TextValue "D22" "72.53"
Set-TextValue "E22" "  -0.59%  "
Set-TextValue "D23" "234.18"
Set-TextValue "E23" "  -0.51%  "
Set-TextValue "D24" "2.21"
Set-TextValue "E24" "  +9.46%  "
Set-TextValue "E25" "  +7.04%  "
Set-TextValue "E26" "  +0.04%  "
Set-TextValue "D27" "2.50"
Set-TextValue "E27" "  +6.88%  "
Set-TextValue "E28" "  +5.10%  "
Set-TextValue "D29" "2.20"
Set-TextValue "E29" "  +1.02%  "
Set-TextValue "D30" "171.22"
Set-TextValue "E30" "  +1.34%  "
Set-TextValue "D31" "20.61"
Set-TextValue "E31" "  +0.45%  "
Set-TextValue "D32" "0.121"
Set-TextValue "E32" "  +2.39%  "
Set-TextValue "E33" "  -1.88%  "
Set-TextValue "D34" "5.54"
Set-TextValue "E34" "  +4.32%  "
Set-TextValue "D35" "0.0719"
Set-TextValue "E35" "  +1.48%  "
Set-TextValue "D36" "4.68"
Set-TextValue "E36" "  -1.26%  "
Set-TextValue "D37" "26.06"
Set-TextValue "E37" "  +19.95%  "
Set-TextValue "D38" "3.93"
Set-TextValue "E38" "  +9.95%  "
Set-TextValue "D39" "0.0295"
Set-TextValue "E39" "  +12.27%  "
Set-TextValue "E40" "  +2.04%  "
Set-TextValue "D41" "68.51"
Set-TextValue "E41" "  +4.37%  "
Set-TextValue "D42" "5.93"
Set-TextValue "E42" "  -0.62%  "
Set-TextValue "D43" "0.209"
Set-TextValue "E43" "  +10.54%  "
Set-TextValue "D44" "11.77"
Set-TextValue "E44" "  +18.87%  "
Set-TextValue "D45" "4.85"
Set-TextValue "E45" "  -2.20%  "
Set-TextValue "D46" "8.79"
Set-TextValue "E46" "  -1.26%  "
Set-TextValue "E47" "  +11.59%  "
Set-TextValue "E48" "  +1.38%  "
Set-TextValue "E49" "  +0.07%  "
Set-TextValue "E50" "  +7.34%  "
Set-TextValue "E51" "  +2.04%  "
